# Insert a new weekly price record as row 15, pushing the existing
# rows 15-61 down to 16-62 (dimension grows from A1:R61 to A1:R62).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 15..61 down by one row, then populate the freed row 15
# with the new weekly record.
$ws.Rows.Item(15).Insert()

$ws.Cells.Item(15, 1).Value = 7
$ws.Cells.Item(15, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(15, 3).Value = "Ñuble"
$ws.Cells.Item(15, 4).Value = 44600
$ws.Cells.Item(15, 4).NumberFormat = $ws.Cells.Item(16, 4).NumberFormat
$ws.Cells.Item(15, 5).Value = 16
$ws.Cells.Item(15, 6).Value = 100112022
$ws.Cells.Item(15, 7).Value = "Arveja Verde"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 60
$ws.Cells.Item(15, 11).Value = 23000
$ws.Cells.Item(15, 12).Value = 24000
$ws.Cells.Item(15, 13).Value = 23500
$ws.Cells.Item(15, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(15, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(15, 16).Value = 940
$ws.Cells.Item(15, 17).Value = 25
$ws.Cells.Item(15, 18).Value = "Hortaliza"
